$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.884.94"
$ws.Range("E2").Value = "  -5.63%  "
$ws.Range("D3").Value = "3.219.98"
$ws.Range("E3").Value = "  -8.85%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.45"
$ws.Range("E5").Value = "  -5.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.63"
$ws.Range("E6").Value = "  -12.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "3.213.21"
$ws.Range("E8").Value = "  -8.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"
$ws.Range("E9").Value = "  -11.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  -12.58%  "
$ws.Range("E11").Value = "  -8.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.499"
$ws.Range("E12").Value = "  -15.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.48"
$ws.Range("E13").Value = "  -17.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000244"
$ws.Range("E14").Value = "  -11.90%  "
$ws.Range("D15").Value = "3.738.35"
$ws.Range("E15").Value = "  -8.90%  "
$ws.Range("D16").Value = "66.796.98"
$ws.Range("E16").Value = "  -5.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "544.67"
$ws.Range("E17").Value = "  -11.61%  "
$ws.Range("D18").Value = "3.215.33"
$ws.Range("E18").Value = "  -9.10%  "
$ws.Range("E19").Value = "  -5.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.14"
$ws.Range("E20").Value = "  -15.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.08"
$ws.Range("E21").Value = "  -15.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.757"
$ws.Range("E22").Value = "  -14.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.75"
$ws.Range("E23").Value = "  -14.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.80"
$ws.Range("E24").Value = "  -12.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.47"
$ws.Range("E25").Value = "  -14.70%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.15"
$ws.Range("E27").Value = "  -16.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.08"
$ws.Range("E28").Value = "  -11.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "29.32"
$ws.Range("E29").Value = "  -13.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.13"
$ws.Range("E30").Value = "  -18.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.58"
$ws.Range("E31").Value = "  -15.25%  "
$ws.Range("E32").Value = "  -13.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "542.36"
$ws.Range("E33").Value = "  -11.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.56"
$ws.Range("E34").Value = "  -19.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.72"
$ws.Range("E35").Value = "  -16.80%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.25"
$ws.Range("E37").Value = "  -6.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0429"
$ws.Range("E38").Value = "  -9.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0846"
$ws.Range("E39").Value = "  -16.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.17"
$ws.Range("E40").Value = "  -15.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.125"
$ws.Range("E41").Value = "  -13.71%  "
$ws.Range("D42").Value = "2.934.32"
$ws.Range("E42").Value = "  -13.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.60"
$ws.Range("E43").Value = "  -26.87%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0586"
$ws.Range("E44").Value = "  -20.96%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.261"
$ws.Range("E45").Value = "  -16.83%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.37"
$ws.Range("E46").Value = "  -19.49%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.98"
$ws.Range("E48").Value = "  -19.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.11"
$ws.Range("E49").Value = "  -18.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.113"
$ws.Range("E50").Value = "  -13.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "123.74"
$ws.Range("E51").Value = "  -7.76%  "
